# Auto-generated edit script to update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text formatting (matches source inlineStr cells)
$cellRefs = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'E7', 'D8', 'E8', 'E9', 'E10', 'E11', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'E18', 'D19', 'E19', 'B20', 'C20', 'D20', 'E20', 'B21', 'C21', 'D21', 'E21', 'E22', 'D23', 'E23', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'E31', 'E32', 'E33', 'D34', 'E34', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'E39', 'E40', 'E41', 'E42', 'D43', 'E43', 'E44', 'E45', 'E46', 'B47', 'C47', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'D49', 'E49', 'E50', 'E51')
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '62.029.73'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.415.46'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '562.97'
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D6').Value = '142.51'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.530'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  +0.55%  '
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  -3.23%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '25.59'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').Value = '0.0000174'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').Value = '2.854.93'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '61.929.48'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '2.413.88'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('E18').Value = '  +1.61%  '
$ws.Range('D19').Value = '321.21'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '4.13'
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.82'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '66.16'
$ws.Range('E23').Value = '  +1.93%  '
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').Value = '8.65'
$ws.Range('E25').Value = '  -4.79%  '
$ws.Range('D26').Value = '564.92'
$ws.Range('E26').Value = '  -2.27%  '
$ws.Range('D27').Value = '2.533.89'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').Value = '0.0₃0933'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').Value = '8.18'
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('E31').Value = '  -3.83%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Value = '1.52'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('D36').Value = '4.77'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').Value = '153.12'
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('D38').Value = '5.44'
$ws.Range('E38').Value = '  -3.77%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').Value = '148.83'
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('E46').Value = '  -2.43%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '19.85'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.594'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('D49').Value = '0.0920'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('E51').Value = '  +0.68%  '
